$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data was recomputed with new TPM values. Sending cluster "ECs" rows
# were dropped, leaving only FAPs/MuSCs as sending clusters (paired with all three
# target clusters). Clear the old rows 2-10 first (drops the stale shared strings too).
$ws.Range("A2:T10").Clear()

# Write the Sending/Ligand/Receptor/Target columns column-by-column so the shared-string
# table is (re)built in the same order as the source export: FAPs, MuSCs, Fbln1, Itgb1, ECs.
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "MuSCs"
$ws.Range("A6").Value = "MuSCs"
$ws.Range("A7").Value = "MuSCs"

$ws.Range("B2").Value = "Fbln1"
$ws.Range("B3").Value = "Fbln1"
$ws.Range("B4").Value = "Fbln1"
$ws.Range("B5").Value = "Fbln1"
$ws.Range("B6").Value = "Fbln1"
$ws.Range("B7").Value = "Fbln1"

$ws.Range("C2").Value = "Itgb1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("C7").Value = "Itgb1"

$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"

# Fill in the numeric measurement columns E:T for each row.
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 72.45391466666668
$ws.Range("H2").Value = 217.361744
$ws.Range("I2").Value = 0.9542169410525404
$ws.Range("J2").Value = 0.9542169410525405
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 4422.703336391594
$ws.Range("R2").Value = 39804.33002752435
$ws.Range("S2").Value = 0.1950050585043228
$ws.Range("T2").Value = 0.1950050585043229

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 72.45391466666668
$ws.Range("H3").Value = 217.361744
$ws.Range("I3").Value = 0.9542169410525404
$ws.Range("J3").Value = 0.9542169410525405
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 7702.899247396235
$ws.Range("R3").Value = 69326.09322656611
$ws.Range("S3").Value = 0.3396348803302162
$ws.Range("T3").Value = 0.3396348803302163

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 72.45391466666668
$ws.Range("H4").Value = 217.361744
$ws.Range("I4").Value = 0.9542169410525404
$ws.Range("J4").Value = 0.9542169410525405
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 9515.98190229307
$ws.Range("R4").Value = 85643.83712063763
$ws.Range("S4").Value = 0.4195770022180013
$ws.Range("T4").Value = 0.4195770022180014

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.476318333333333
$ws.Range("H5").Value = 10.428955
$ws.Range("I5").Value = 0.04578305894745947
$ws.Range("J5").Value = 0.04578305894745948
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 212.2000551926828
$ws.Range("R5").Value = 1909.800496734145
$ws.Range("S5").Value = 0.009356287553130555
$ws.Range("T5").Value = 0.009356287553130557

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.476318333333333
$ws.Range("H6").Value = 10.428955
$ws.Range("I6").Value = 0.04578305894745947
$ws.Range("J6").Value = 0.04578305894745948
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("Q6").Value = 369.5829272543433
$ws.Range("R6").Value = 3326.24634528909
$ws.Range("S6").Value = 0.01629558549822001
$ws.Range("T6").Value = 0.01629558549822002

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.476318333333333
$ws.Range("H7").Value = 10.428955
$ws.Range("I7").Value = 0.04578305894745947
$ws.Range("J7").Value = 0.04578305894745948
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 456.5741202363044
$ws.Range("R7").Value = 4109.16708212674
$ws.Range("S7").Value = 0.02013118589610891
$ws.Range("T7").Value = 0.02013118589610891
